$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the duplicated rows (rows 5-11 repeat "Museu Histórico Nacional/Ibram"
# with the same license text already present on row 4). Deleting them shifts
# all subsequent rows up by 7, matching the new A1:B15 data range.
$ws.Rows("5:11").Delete()

# Reset the view: scroll back to the top and select the full used range.
$ws.Range("A1:B15").Select()
